$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 260 (existing rows 260-279 shift down to 262-281)
$ws.Rows.Item(260).Insert()
$ws.Rows.Item(260).Insert()

# Copy formatting (incl. date number format on column D) from the row that is
# now at 262 (the old row 260) down onto the two freshly inserted rows.
$ws.Range("A262:R262").Copy()
$ws.Range("A260:R260").PasteSpecial(-4122)
$ws.Range("A262:R262").Copy()
$ws.Range("A261:R261").PasteSpecial(-4122)

# Row 260 - new "Primera" entry
$ws.Cells.Item(260, 1).Value = 3
$ws.Cells.Item(260, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(260, 3).Value = "Coquimbo"
$ws.Cells.Item(260, 4).Value = 44461
$ws.Cells.Item(260, 5).Value = 5
$ws.Cells.Item(260, 6).Value = 100112037
$ws.Cells.Item(260, 7).Value = "Cebollín"
$ws.Cells.Item(260, 8).Value = "Sin especificar"
$ws.Cells.Item(260, 9).Value = "Primera"
$ws.Cells.Item(260, 10).Value = 120
$ws.Cells.Item(260, 11).Value = 3500
$ws.Cells.Item(260, 12).Value = 3500
$ws.Cells.Item(260, 13).Value = 3500
$ws.Cells.Item(260, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(260, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(260, 16).Value = 97
$ws.Cells.Item(260, 17).Value = 36
$ws.Cells.Item(260, 18).Value = "Hortaliza"

# Row 261 - new "Segunda" entry
$ws.Cells.Item(261, 1).Value = 3
$ws.Cells.Item(261, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(261, 3).Value = "Coquimbo"
$ws.Cells.Item(261, 4).Value = 44461
$ws.Cells.Item(261, 5).Value = 5
$ws.Cells.Item(261, 6).Value = 100112037
$ws.Cells.Item(261, 7).Value = "Cebollín"
$ws.Cells.Item(261, 8).Value = "Sin especificar"
$ws.Cells.Item(261, 9).Value = "Segunda"
$ws.Cells.Item(261, 10).Value = 160
$ws.Cells.Item(261, 11).Value = 2500
$ws.Cells.Item(261, 12).Value = 2500
$ws.Cells.Item(261, 13).Value = 2500
$ws.Cells.Item(261, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(261, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(261, 16).Value = 69
$ws.Cells.Item(261, 17).Value = 36
$ws.Cells.Item(261, 18).Value = "Hortaliza"

Write-Output "done"
